$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct existing values ---
$ws.Range("D2").Value = 28.1
$ws.Range("B3").Value = 94.59999999999999
$ws.Range("D3").Value = 28.1

# Fill in the previously-missing evening readings for 2025-11-10 (row 5)
$ws.Range("C5").Value = 95.3
$ws.Range("E5").Value = 26.8

# --- Append new rows 6-8 ---
# Copy the formatting/layout of row 5 down into rows 6:8 first, so the new
# rows inherit the date style (A) and the blank placeholder cells (F:H).
$ws.Range("A5:H5").Copy($ws.Range("A6:A8"))

# Row 6 - 2025-11-11
$ws.Range("A6").Value = 45972
$ws.Range("B6").Value = 94.40000000000001
$ws.Range("C6").Value = 95.2
$ws.Range("D6").Value = 27.9
$ws.Range("E6").Value = 27.1

# Row 7 - 2025-11-12
$ws.Range("A7").Value = 45973
$ws.Range("B7").Value = 94.40000000000001
$ws.Range("C7").Value = 95.09999999999999
$ws.Range("D7").Value = 27.8
$ws.Range("E7").Value = 27.7

# Row 8 - 2025-11-13
$ws.Range("A8").Value = 45974
$ws.Range("B8").Value = 94.40000000000001
$ws.Range("C8").Value = 95.8
$ws.Range("D8").Value = 28
$ws.Range("E8").Value = 26.9
